# Lac_AllBounds14.xlsx — "Added Flow vs R1L to the cell data modeled by tissue slice code"
#
# This mirrors the existing Kpl summary block (rows 22-24, cols B-F) but for
# the Flow_Lac series (column F), landing in rows 38-40, cols F-J.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- best-effort VBA codenames (workbookPr/sheetPr codeName) ---------------
# These match the default names Excel assigns a fresh single-sheet workbook,
# so setting them explicitly is a harmless no-op if the host doesn't persist
# them, and correct if it does.
try { $wb.CodeName = "ThisWorkbook" } catch {}
try { $ws.CodeName = "Sheet1" } catch {}

# --- header row (38): tissue-slice / cell-line labels -----------------------
$ws.Range("G38").Value = "HK-2"
$ws.Range("H38").Value = "UMRC6"
$ws.Range("I38").Value = "UOK262"
$ws.Range("J38").Value = "UOK + DIDS"

# --- row 39: series label + AVERAGE of Flow_Lac (column F) per group -------
$ws.Range("F39").Value = "Flow_Lac"
$ws.Range("G39").Formula = "=AVERAGE(F`$1:F`$3)"
$ws.Range("H39").Formula = "=AVERAGE(F`$4:F`$6)"
$ws.Range("I39").Formula = "=AVERAGE(F`$9:F`$11)"
$ws.Range("J39").Formula = "=AVERAGE(F`$13:F`$16)"

# --- row 40: standard error of the mean per group ---------------------------
$ws.Range("G40").Formula = "=STDEV(F`$1:F`$3)/SQRT(COUNT(F`$1:F`$3))"
$ws.Range("H40").Formula = "=STDEV(F`$4:F`$6)/SQRT(COUNT(F`$4:F`$6))"
$ws.Range("I40").Formula = "=STDEV(F`$9:F`$11)/SQRT(COUNT(F`$9:F`$11))"
$ws.Range("J40").Formula = "=STDEV(F`$13:F`$16)/SQRT(COUNT(F`$13:F`$16))"

# --- view state: scroll + selection to the new block ------------------------
try {
    $excel.ActiveWindow.ScrollRow = 19
    $excel.ActiveWindow.ScrollColumn = 1
} catch {}

$ws.Range("F38:J40").Select()
